$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "248.84"
    "D4"  = "5.375"
    "D5"  = "0.05682"
    "D6"  = "3.402"
    "D7"  = "6.339"
    "D8"  = "0.8060"
    "D9"  = "0.9168"
    "D11" = "0.07448"
    "D12" = "0.03110"
    "D13" = "0.03049"
    "D14" = "0.09382"
    "D17" = "0.04787"
    "D19" = "0.0005852"
    "D20" = "0.006497"
    "D21" = "0.004987"
    "D22" = "0.0009997"
    "D23" = "0.0001501"
    "D24" = "3.696"
    "D26" = "0.3253"
    "D27" = "0.1307"
    "D40" = "0.04004"
    "D41" = "0.1071"
    "D42" = "0.002741"
    "D43" = "0.003041"
    "D44" = "0.007971"
    "D45" = "0.00005753"
    "D47" = "0.4992"
    "D48" = "0.2109"
    "D49" = "0.00002101"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
